$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 18, pushing existing rows 18-22 down to 19-23
$ws.Rows.Item(18).Insert()

# Fill the new row 18 with the weekly data point (same categorical data as the
# row that used to be there, but with an updated date and prices)
$ws.Cells.Item(18, 1).Value = 4
$ws.Cells.Item(18, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(18, 3).Value = "Los Lagos"
$ws.Cells.Item(18, 4).Value = 44726
$ws.Cells.Item(18, 5).Value = 10
$ws.Cells.Item(18, 6).Value = 100112012
$ws.Cells.Item(18, 7).Value = "Espinaca"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 30
$ws.Cells.Item(18, 11).Value = 14000
$ws.Cells.Item(18, 12).Value = 14000
$ws.Cells.Item(18, 13).Value = 14000
$ws.Cells.Item(18, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(18, 15).Value = "Región Metropolitana"
$ws.Cells.Item(18, 16).Value = 1400
$ws.Cells.Item(18, 17).Value = 10
$ws.Cells.Item(18, 18).Value = "Hortaliza"

# Match the date formatting/style used by the other date cells in column D
$ws.Cells.Item(18, 4).NumberFormat = $ws.Cells.Item(19, 4).NumberFormat
